# Generate Report for Handback
# Updates the localization-status workbook to reflect that the handback
# has completed: status text changes, target/handback file + datetime
# columns get populated (with hyperlinks) on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Update the "Status" text everywhere it appears (Overview E/F cols,
#    and column C on each language sheet) from "Ready for handoff" to
#    the new handback status.
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Widen the Status-related columns now that the text is longer.
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14
$wsZhCn.Columns.Item(3).ColumnWidth = 29.14
$wsDeDe.Columns.Item(3).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# 2. Populate "Latest Target File" (I) / "Latest Handback File" (J) /
#    "Latest Handback DateTime" (K) for both rows of each language
#    sheet.
# ---------------------------------------------------------------------
function Set-HandbackRow($ws, $row, $targetFileName, $handbackFile, $handbackDateTime) {
    $ws.Cells.Item($row, 9).Value = $targetFileName        # I - Latest Target File
    $ws.Cells.Item($row, 10).Value = $handbackFile         # J - Latest Handback File
    $ws.Cells.Item($row, 11).Value = $handbackDateTime     # K - Latest Handback DateTime
}

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a277dfe0a1bd0995aa5b82a4d82d49aa98a74c39/e2e/"

$file1Name = "4d7687bb-c846-4551-851a-22c87b354c35.md"
$file2Name = "befad8eb-c126-4f04-bb13-ab104788cedd.md"
$file1Url = $repoBase + $file1Name
$file2Url = $repoBase + $file2Name

# zh-cn sheet
Set-HandbackRow $wsZhCn 2 $file1Name "4d7687bb-c846-4551-851a-22c87b354c35.daa8c95460913e4befa3f5ac44593de6ed74ad7b.zh-cn.xlf" "2016-08-26 04:46:13"
Set-HandbackRow $wsZhCn 3 $file2Name "befad8eb-c126-4f04-bb13-ab104788cedd.1a66f051173d7a02929dd191424535fc6fab8bed.zh-cn.xlf" "2016-08-26 04:46:13"

# de-de sheet
Set-HandbackRow $wsDeDe 2 $file1Name "4d7687bb-c846-4551-851a-22c87b354c35.daa8c95460913e4befa3f5ac44593de6ed74ad7b.de-de.xlf" "2016-08-26 04:46:20"
Set-HandbackRow $wsDeDe 3 $file2Name "befad8eb-c126-4f04-bb13-ab104788cedd.1a66f051173d7a02929dd191424535fc6fab8bed.de-de.xlf" "2016-08-26 04:46:20"

# Columns I (Target File) and J (Handback File) now hold file names, so
# give them the same fixed width used for the other file-name columns.
$wsZhCn.Columns.Item(9).ColumnWidth = 39.16
$wsZhCn.Columns.Item(10).ColumnWidth = 39.16
$wsDeDe.Columns.Item(9).ColumnWidth = 39.16
$wsDeDe.Columns.Item(10).ColumnWidth = 39.16

# ---------------------------------------------------------------------
# 3. Rebuild the hyperlinks on each language sheet: Source File Name
#    (A) already linked to the source .md file on GitHub; now Latest
#    Target File (I) gets the very same link, row by row.
# ---------------------------------------------------------------------
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $file1Url, $null, $null, $file1Name)
    $ws.Hyperlinks.Add($ws.Range("I2"), $file1Url, $null, $null, $file1Name)
    $ws.Hyperlinks.Add($ws.Range("A3"), $file2Url, $null, $null, $file2Name)
    $ws.Hyperlinks.Add($ws.Range("I3"), $file2Url, $null, $null, $file2Name)

    # Match the look of the existing hyperlinked cells (underline + the
    # workbook's custom hyperlink blue) on the newly linked cells.
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 0xED9564
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = 0xED9564
}
